# Updated Today's Task in TasksBreakDown/AashishSharma.xlsx
# Fill in "Hours Burnt" (column F) for the tasks that were worked on today.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> Hours Burnt value to enter (column F)
$updates = @{
    "F3"  = 1
    "F4"  = 1
    "F5"  = 1
    "F21" = 1
    "F22" = 1
    "F23" = 1
    "F25" = 1
    "F26" = 2
    "F28" = 1
    "F29" = 1
    "F30" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Update the view state to match where the user was working: scrolled down
# so row 14 is at the top, with F22 as the active selected cell.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F22").Select()
